# VGA Control architecture and testbench
# Add a new time-record row (row 16) for "30.3.2020 / VGA Control / RTL and TB",
# following the same pattern/format as the existing rows 11-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 15's formatting into a new row 16 (keeps styles/number formats
# identical to the existing block instead of Excel re-minting new style indices).
$ws.Rows("15:15").Copy()
$ws.Rows("16:16").Insert(-4121)  # xlShiftDown
$excel.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("A16").Value = "30.3.2020"
$ws.Range("B16").Value = 0.58333333333333337
$ws.Range("C16").Value = 0.6875
$ws.Range("D16").Formula = "=C16-B16"
$ws.Range("E16").Value = "VGA Control"
$ws.Range("F16").Value = "RTL and TB"

# Move the active selection down to the next empty row, like in the source file.
$ws.Range("A17").Select()
